$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new data row at row 2, shifting every existing row (2-25)
# down by one (to 3-26) while keeping each row's own formatting intact.
$ws.Rows(2).Insert()

# The freshly inserted row inherited the header row's bold/centered style;
# strip that back to the plain/unstyled look used by the other data rows.
$ws.Range("A2:R2").ClearFormats()

# Column D holds dates, formatted the same way as the rest of the date
# column.
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row with the latest weekly observation.
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "Vega Monumental Concepción"
$ws.Range("C2").Value = "Bíobío"
$ws.Range("D2").Value = 44631
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 100112022
$ws.Range("G2").Value = "Arveja Verde"
$ws.Range("H2").Value = "Perfection"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 24000
$ws.Range("L2").Value = 25000
$ws.Range("M2").Value = 24467
$ws.Range("N2").Value = "`$/saco 25 kilos"
$ws.Range("O2").Value = "Carahue"
$ws.Range("P2").Value = 979
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = "Hortaliza"
